$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor values update
$ws.Range("B3").Value = 0.0124905995992804
$ws.Range("C3").Value = 0.01277285751795381
$ws.Range("D3").Value = 0.051592380313197

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values update
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.014213000018349
$ws.Range("C4").Value = 0.01454834688729477
$ws.Range("D4").Value = 0.05998628857852755

# Row 5: AdaBoostRegressor -> MLPRegressor, values update
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01298542539928405
$ws.Range("C5").Value = 0.01292917550587768
$ws.Range("D5").Value = 0.01867569471876836
